$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "72.925.82"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.03%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "4.048.36"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  -0.18%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "561.70"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +4.83%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "151.25"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "4.043.29"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.73%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.697"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("E9").Value = "  -0.08%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.766"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.172"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "53.88"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +13.53%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000328"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.89%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "10.97"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.43%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.698.17"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "4.056.08"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "14.53"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.19%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "20.83"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("E19").Value = "  +1.50%  "

$ws.Range("E20").Value = "  -0.41%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "72.899.36"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.15%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "446.08"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +3.98%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "98.30"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.02%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.55"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "4.42"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.86%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "14.82"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "4.38"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +21.11%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "11.37"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "10.99"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.13%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.95"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.54%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "37.36"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.94"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +14.24%  "

$ws.Range("E33").Value = "  +3.52%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "13.65"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.29%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "689.11"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "49.07"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +14.80%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "67.43"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.65%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0922"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +11.96%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.452"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +6.21%  "

$ws.Range("E40").Value = "  -2.17%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.43"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.61%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.40"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.02%  "

$ws.Range("B43").Value = "Dai"
$ws.Range("C43").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0498"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +2.00%  "

$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "11.04"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +14.70%  "

$ws.Range("E46").Value = "  +0.10%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.152"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("E48").Value = "  +4.89%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.57"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +5.94%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "3.14"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +4.70%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "3.32"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.88%  "
